$wb = $excel.ActiveWorkbook

$wsPatient = $wb.Worksheets.Item("patient")
$wsCellModel = $wb.Worksheets.Item("cell_model")
$wsSharing = $wb.Worksheets.Item("sharing")

# --- cell_model sheet: update model_id values in A2/A3 ---
# The two existing model_id cells originally used a style (s=24) that is
# distinct from the normal body-cell style (s=4) used everywhere else in
# the workbook. Re-typing the values also brings the formatting in line
# with the rest of the workbook, so pull the common body style across
# first (copy format only) and then overwrite the values.
$wsPatient.Range("B3").Copy()
$wsCellModel.Range("A2:A3").PasteSpecial(-4122)

$wsCellModel.Range("A2").Value = "CRC0228PR"
$wsCellModel.Range("A3").Value = "CRC0228PRaS"

# --- sharing sheet: update model_id values in A2/A3 ---
$wsSharing.Range("A2").Value = "CRC0228PR"
$wsSharing.Range("A3").Value = "CRC0228PRaS"

# --- selections / view state ---
$wsSharing.Range("A2").Select()

$wsCellModel.Activate()
$wsCellModel.Range("A2").Select()
